$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the values in row 5 and row 6 for columns A, B, E, F, G, H
$cols = @("A", "B", "E", "F", "G", "H")

foreach ($col in $cols) {
    $range5 = $ws.Range($col + "5")
    $range6 = $ws.Range($col + "6")
    $tmp = $range5.Value2
    $range5.Value2 = $range6.Value2
    $range6.Value2 = $tmp
}
